$wb = $excel.ActiveWorkbook

$annSheet  = $wb.Worksheets.Item("ANN")
$svmSheet  = $wb.Worksheets.Item("SVM")
$ksvmSheet = $wb.Worksheets.Item("K-SVMeans")

# --- K-SVMeans sheet ("No Volume" results added in column F) ---
$ksvmSheet.Range("F4").Value = 59.42
$ksvmSheet.Range("F5").Value = 55.67
$ksvmSheet.Range("F6").Value = 51.87
$ksvmSheet.Range("F7").Value = 57.61
$ksvmSheet.Range("F8").Value = 55.84
$ksvmSheet.Range("F9").Formula = "=AVERAGE(F4:F8)"

# F6 carries the same "highlighted" cell format already used on E8 (style s="4")
$ksvmSheet.Range("E8").Copy() | Out-Null
$ksvmSheet.Range("F6").PasteSpecial(-4122) | Out-Null

# F9 carries the same totals-row format already used on E9 (style s="2")
$ksvmSheet.Range("E9").Copy() | Out-Null
$ksvmSheet.Range("F9").PasteSpecial(-4122) | Out-Null
$ksvmSheet.Range("F9").Formula = "=AVERAGE(F4:F8)"

$ksvmSheet.Application.CutCopyMode = $false

# K-SVMeans sheet gets its own explicit page setup (portrait), as added in the diff
$ksvmSheet.PageSetup.Orientation = 1

# --- Selections on each sheet ---
$annSheet.Range("H9").Select() | Out-Null
$svmSheet.Range("G6").Select() | Out-Null

# K-SVMeans becomes the active tab/sheet (selected last)
$ksvmSheet.Range("H11").Select() | Out-Null

$wb.Save() | Out-Null
